$wb = $excel.ActiveWorkbook

# Add a new worksheet named "Instructions" after Sheet1
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Instructions"

# Fill in the instructions content. The write order below matches the
# order the strings were first authored (not plain row order) so the
# resulting shared-string table is built up in the same sequence as the
# target workbook. Leading "'" forces text so Excel doesn't coerce the
# step numbers like "1." into numeric values; each such cell's style is
# immediately reset to Normal afterwards so it doesn't keep Excel's
# "quote prefix" styling (restores a plain, unstyled text cell).
$ws2.Range("A1").Value = "In order to go from excel to html you need to :"

$ws2.Range("A2").Value = "'1."
$ws2.Range("A2").Style = "Normal"
$ws2.Range("B2").Value = "Copy the excel cells to MS word"

$ws2.Range("A3").Value = "'2."
$ws2.Range("A3").Style = "Normal"
$ws2.Range("B4").Value = "Save the MS word as html"

$ws2.Range("A4").Value = "'3."
$ws2.Range("A4").Style = "Normal"
$ws2.Range("B5").Value = "Open the html in notepad++"

$ws2.Range("A5").Value = "'4. "
$ws2.Range("A5").Style = "Normal"
$ws2.Range("B7").Value = "Replace all \n by nothing"

$ws2.Range("A6").Value = "'5. "
$ws2.Range("A6").Style = "Normal"
$ws2.Range("B8").Value = "Replace all \r by nothing"

$ws2.Range("A7").Value = "'6. "
$ws2.Range("A7").Style = "Normal"
$ws2.Range("B9").Value = "Use pretty XML formatting to format the table"

$ws2.Range("B6").Value = "Copy the table element in a new file"

$ws2.Range("A8").Value = "'7."
$ws2.Range("A8").Style = "Normal"

$ws2.Range("A9").Value = "'8."
$ws2.Range("A9").Style = "Normal"
$ws2.Range("B10").Value = "In table element, replace the style border-collapse: collapse by 'border-collapse: unset'"

$ws2.Range("B3").Value = "Make sure the cells are well displayed (readjust the with of cells if needed)"

$ws2.Range("A10").Value = "'9."
$ws2.Range("A10").Style = "Normal"

$ws2.Range("A11").Value = "'10."
$ws2.Range("A11").Style = "Normal"
$ws2.Range("B12").Value = "Replace 'background:#D9E1F2;' by 'background:#D9E1F2;vertical-align: middle;'"

$ws2.Range("A12").Value = "'11."
$ws2.Range("A12").Style = "Normal"
$ws2.Range("B13").Value = "Replace 'background:#C6E0B4;' by 'background:#C6E0B4;vertical-align: middle;'"

$ws2.Range("A13").Value = "'12."
$ws2.Range("A13").Style = "Normal"
$ws2.Range("B11").Value = "Remove from table element the style: margin-left:21.0pt;"

# Set selections to match target state
$ws1.Range("W12").Select()
$ws2.Range("B12").Select()

# Make Instructions sheet the active sheet/tab
$ws2.Activate()
